$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 07:04"

# Update country rows: some countries overtook their neighbors in total-case ranking,
# so their row (country name + stats) moved up, pushing others down unchanged.

# Row 15: Pakistan
$ws.Range("A15").Value = "Pakistan"
$ws.Range("B15").Value = 277402
$ws.Range("C15").Value = 1114
$ws.Range("D15").Value = 246131
$ws.Range("E15").Value = 25347
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 32
$ws.Range("H15").Value = 5924

# Row 55: Kirguistan
$ws.Range("A55").Value = "Kirguistan"
$ws.Range("B55").Value = 35223
$ws.Range("C55").Value = 631
$ws.Range("D55").Value = 23985
$ws.Range("E55").Value = 9874
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 17
$ws.Range("H55").Value = 1364

# Row 56: Ghana
$ws.Range("A56").Value = "Ghana"
$ws.Range("B56").Value = 35142
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 31286
$ws.Range("E56").Value = 3681
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 175

# Row 57: Suiza
$ws.Range("A57").Value = "Suiza"
$ws.Range("B57").Value = 34802
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 31100
$ws.Range("E57").Value = 1723
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 1979

# Row 70: Venezuela
$ws.Range("A70").Value = "Venezuela"
$ws.Range("B70").Value = 17158
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 10421
$ws.Range("E70").Value = 6581
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 156

# Row 71: Costa Rica
$ws.Range("A71").Value = "Costa Rica"
$ws.Range("B71").Value = 16800
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 4050
$ws.Range("E71").Value = 12617
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 133

# Row 109: Tailandia
$ws.Range("A109").Value = "Tailandia"
$ws.Range("B109").Value = 3304
$ws.Range("C109").Value = 6
$ws.Range("D109").Value = 3111
$ws.Range("E109").Value = 135
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 58

# Row 187: Barbados
$ws.Range("A187").Value = "Barbados"
$ws.Range("B187").Value = 110
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 95
$ws.Range("E187").Value = 8
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 7

# Row 188: Islas Turcas y Caicos
$ws.Range("A188").Value = "Islas Turcas y Caicos"
$ws.Range("B188").Value = 104
$ws.Range("C188").Value = 5
$ws.Range("D188").Value = 37
$ws.Range("E188").Value = 65
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 2

# Row 189: Butan
$ws.Range("A189").Value = "Butan"
$ws.Range("B189").Value = 101
$ws.Range("C189").Value = 2
$ws.Range("D189").Value = 88
$ws.Range("E189").Value = 13
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 0

# Row 192: Papua Nueva Guinea
$ws.Range("A192").Value = "Papua Nueva Guinea"
$ws.Range("B192").Value = 63
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 11
$ws.Range("E192").Value = 50
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 2

# Row 203: Granada
$ws.Range("A203").Value = "Granada"
$ws.Range("B203").Value = 24
$ws.Range("C203").Value = 1
$ws.Range("D203").Value = 23
$ws.Range("E203").Value = 1
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 204: Timor Oriental
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("B204").Value = 24
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 24
$ws.Range("E204").Value = 0
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0
